# COVID-19-TW-PlaceCode.xlsx update
# Appends 5 more days of data (2020-04-28 .. 2020-05-02) to the "CodeBook"
# sheet, extending the table from column AI to column AL, and backfills the
# already-existing-but-unpopulated AH/AI placeholder columns with real data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CodeBook")

$newCols = @("AH", "AI", "AJ", "AK", "AL")

# --- 1. Copy the formatting of the last fully-populated column (AG) onto
#        the five target columns so every cell picks up the same style
#        (date style on row 1, header style on row 2, data style on rows
#        3-24, SUM-formula style on row 25) without minting new styles.
$ws.Range("AG1:AG25").Copy()
$ws.Range("AH1:AL25").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- 2. Row 1: the five new date headers (serial dates 2020-04-28 .. 2020-05-02)
$newDates = @{
    "AH" = 43949
    "AI" = 43950
    "AJ" = 43951
    "AK" = 43952
    "AL" = 43953
}
foreach ($col in $newCols) {
    $ws.Range($col + "1").Value = $newDates[$col]
}

# --- 3. Row 2: repeat the "病例數" (case count) header label
foreach ($col in $newCols) {
    $ws.Range($col + "2").Value = "病例數"
}

# --- 4. Rows 3-24: per-county case counts for each of the 5 new days.
#        Most days simply repeat the last known total; a couple of counties
#        saw their counts tick up on the final day (2020-05-02).
$rowData = @{
    3  = @(116, 116, 116, 116, 116)   # Taipei City
    4  = @(40,  40,  40,  40,  40)    # New Taipei City
    5  = @(16,  16,  16,  16,  16)    # Taichung City
    6  = @(44,  44,  44,  44,  46)    # Tainan City
    7  = @(7,   7,   7,   7,   7)     # Kaohsiung City
    8  = @(11,  11,  11,  11,  11)    # Keelung City
    9  = @(4,   4,   4,   4,   4)     # Hsinchu City
    10 = @(90,  90,  90,  90,  90)    # Chiayi City
    11 = @(51,  51,  51,  51,  52)    # New Taipei (county-level row)
    12 = @(6,   6,   6,   6,   6)     # Taoyuan City
    13 = @(2,   2,   2,   2,   2)     # Hsinchu County
    14 = @(3,   3,   3,   3,   3)     # Yilan County
    15 = @(18,  18,  18,  18,  18)    # Miaoli County
    16 = @(2,   2,   2,   2,   2)     # Changhua County
    17 = @(5,   5,   5,   5,   5)     # Nantou County
    18 = @(2,   2,   2,   2,   2)     # Yunlin County
    19 = @(12,  12,  12,  12,  12)    # Chiayi County
    20 = @(0,   0,   0,   0,   0)     # Pingtung County
    21 = @(0,   0,   0,   0,   0)     # Penghu County
    22 = @(0,   0,   0,   0,   0)     # Hualien County
    23 = @(0,   0,   0,   0,   0)     # Taitung County
    24 = @(0,   0,   0,   0,   0)     # Kinmen County
}

foreach ($row in $rowData.Keys) {
    $vals = $rowData[$row]
    for ($i = 0; $i -lt $newCols.Length; $i++) {
        $ws.Range($newCols[$i] + $row).Value = $vals[$i]
    }
}

# --- 5. Row 25: grand-total SUM formulas, same pattern as the rest of the row.
foreach ($col in $newCols) {
    $ws.Range($col + "25").Formula = "=SUM(" + $col + "3:" + $col + "24)"
}

# --- 6. Housekeeping to mirror the widened view/selection Excel leaves
#        behind after scrolling right to work on the newly-added columns.
$ws.Range("AG25:AL25").Select()
$ws.Application.ActiveWindow.ScrollColumn = 26
